$d = $word.ActiveDocument

# First run: trim "5 - The western half ... to wildfires " down to "5 - "
$d.Content.Find.Execute(
    "5 – The western half of the contiguous United States is the most prone to wildfires ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "5 – ",
    2
)

# Second run: replace the old continuation with the new sentence(s)
$d.Content.Find.Execute(
    "due to the effects of climate change on its large stretches of undeveloped land. For those reasons, I will focus on a few case studies on this area. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "When viewed back-to-back in this fashion, drought and increasing fire sizes become a clearly regional issue. These trends converge on the western half of the United States. The importance of climatic changes and their effects on wildfires becomes more evident when you isolate the largest wildfire of each decade, for the past century. ",
    2
)
